# momentum, backprop, CV, FCT
# Added in backprop, momentum, and the hidden layer. CV for free parameters.
# Forced choice test cases for ortho to phono.
#
# Appends a new data row ("diphthong ") to the bottom of the table on
# Sheet1, mirroring the zero-filled rows already present, and updates the
# sheet's scroll position / selection the way Excel does after typing a
# new row at the bottom of a list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 13

# Row label in column A.
$ws.Cells.Item($newRow, 1).Value = "diphthong "

# Feature columns B:Y are all 0 for this new phoneme row.
for ($col = 2; $col -le 25; $col++) {
    $ws.Cells.Item($newRow, $col).Value = 0
}

# Scroll the window down a bit (Excel records this as topLeftCell on the
# sheetView) and leave the freshly-entered row's data cells selected, same
# as what happens after typing across a row and pressing Enter/Tab.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B13:Y13").Select()
